$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each cell is explicitly forced to Text format before assignment, then
# restored to the default "Normal" style, so that numeric-looking strings
# (e.g. "0.9985") are stored as literal text -- matching the original
# inlineStr cells -- instead of being auto-coerced to floating point by Excel.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "30.288.77"
Set-TextValue "E2" "  -0.18%  "
Set-TextValue "D3" "1.931.22"
Set-TextValue "E3" "  -0.51%  "
Set-TextValue "D4" "0.9985"
Set-TextValue "E4" "  -0.25%  "
Set-TextValue "D5" "0.7512"
Set-TextValue "E5" "  +4.16%  "
Set-TextValue "D6" "250.10"
Set-TextValue "E6" "  -0.33%  "
Set-TextValue "D7" "0.9987"
Set-TextValue "E7" "  -0.18%  "
Set-TextValue "B8" "Solana"
Set-TextValue "C8" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D8" "28.14"
Set-TextValue "E8" "  -2.91%  "
Set-TextValue "B9" "Cardano"
Set-TextValue "C9" "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue "D9" "0.3229"
Set-TextValue "E9" "  -3.14%  "
Set-TextValue "D10" "0.07131"
Set-TextValue "E10" "  -3.35%  "
Set-TextValue "D11" "0.7907"
Set-TextValue "E11" "  -3.02%  "
Set-TextValue "D12" "0.08015"
Set-TextValue "E12" "  -1.34%  "
Set-TextValue "D13" "1.934.13"
Set-TextValue "E13" "  -0.30%  "
Set-TextValue "D14" "5.399"
Set-TextValue "E14" "  -1.64%  "
Set-TextValue "D15" "94.56"
Set-TextValue "E15" "  -0.55%  "
Set-TextValue "D16" "14.54"
Set-TextValue "E16" "  -2.80%  "
Set-TextValue "D17" "30.301.90"
Set-TextValue "E17" "  -0.19%  "
Set-TextValue "D18" "253.43"
Set-TextValue "E18" "  +0.37%  "
Set-TextValue "D19" "0.000008074"
Set-TextValue "E19" "  -3.32%  "
Set-TextValue "E20" "  -2.01%  "
Set-TextValue "D21" "2.188.59"
Set-TextValue "E21" "  -0.30%  "
Set-TextValue "D22" "0.9991"
Set-TextValue "E22" "  -0.14%  "
Set-TextValue "D23" "0.9983"
Set-TextValue "E23" "  -0.36%  "
Set-TextValue "D24" "6.844"
Set-TextValue "E24" "  -1.83%  "
Set-TextValue "D25" "9.592"
Set-TextValue "E25" "  -2.63%  "
Set-TextValue "D26" "164.16"
Set-TextValue "B27" "Stellar"
Set-TextValue "C27" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D27" "0.1345"
Set-TextValue "E27" "  +1.92%  "
Set-TextValue "B28" "EthereumClassic"
Set-TextValue "C28" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D28" "19.10"
Set-TextValue "E28" "  -1.61%  "
Set-TextValue "D29" "2.301"
Set-TextValue "E29" "  -4.41%  "
Set-TextValue "D30" "1.354"
Set-TextValue "E30" "  +0.49%  "
Set-TextValue "D31" "1.534"
Set-TextValue "E31" "  -2.46%  "
Set-TextValue "D32" "4.423"
Set-TextValue "E32" "  -0.83%  "
Set-TextValue "D33" "4.157"
Set-TextValue "E33" "  -2.19%  "
Set-TextValue "D34" "1.302"
Set-TextValue "E34" "  -0.19%  "
Set-TextValue "D35" "0.05120"
Set-TextValue "E35" "  -2.87%  "
Set-TextValue "D36" "0.7490"
Set-TextValue "E36" "  -0.92%  "
Set-TextValue "D37" "2.768"
Set-TextValue "E37" "  +0.64%  "
Set-TextValue "D38" "0.01978"
Set-TextValue "E38" "  -0.79%  "
Set-TextValue "D39" "2.799"
Set-TextValue "E39" "  -2.10%  "
Set-TextValue "D40" "78.26"
Set-TextValue "E40" "  -3.94%  "
Set-TextValue "D41" "6.409"
Set-TextValue "E41" "  -3.06%  "
Set-TextValue "D42" "0.4507"
Set-TextValue "E42" "  -1.11%  "
Set-TextValue "D43" "1.996"
Set-TextValue "E43" "  -2.27%  "
Set-TextValue "D44" "0.8442"
Set-TextValue "E44" "  -0.41%  "
Set-TextValue "D45" "0.9991"
Set-TextValue "D46" "101.44"
Set-TextValue "B47" "Maker"
Set-TextValue "C47" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D47" "1.007.01"
Set-TextValue "E47" "  +13.76%  "
Set-TextValue "B48" "EnergySwap"
Set-TextValue "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "9.819"
Set-TextValue "E48" "  -0.26%  "
Set-TextValue "B49" "Aptos"
Set-TextValue "C49" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D49" "7.534"
Set-TextValue "E49" "  +0.34%  "
Set-TextValue "D50" "37.37"
Set-TextValue "E50" "  +1.30%  "
Set-TextValue "D51" "0.06067"
Set-TextValue "E51" "  +0.30%  "
